# Update Fgf7-Fgfr2 sheet with new TPM-based NATMI output values.
# Adds a new target/receptor-expressing cluster "Resolving-Mac" and
# refreshes the specificity/weight statistics across the 3x4 sender/target grid.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: ECs -> ECs
$ws.Range("A2").Value = "ECs"
$ws.Range("B2").Value = "Fgf7"
$ws.Range("C2").Value = "Fgfr2"
$ws.Range("D2").Value = "ECs"
$ws.Range("E2").Value = 1
$ws.Range("F2").Value = 0.3333333333333333
$ws.Range("G2").Value = 0.08805033333333334
$ws.Range("H2").Value = 0.264151
$ws.Range("I2").Value = 0.005589762818257384
$ws.Range("J2").Value = 0.005589762818257385
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 0.714474
$ws.Range("N2").Value = 2.143422
$ws.Range("O2").Value = 0.138796410342318
$ws.Range("P2").Value = 0.138796410342318
$ws.Range("Q2").Value = 0.062909673858
$ws.Range("R2").Value = 0.5661870647220001
$ws.Range("S2").Value = 0.0007758390138390839
$ws.Range("T2").Value = 0.000775839013839084

# Row 3: ECs -> FAPs
$ws.Range("A3").Value = "ECs"
$ws.Range("B3").Value = "Fgf7"
$ws.Range("C3").Value = "Fgfr2"
$ws.Range("D3").Value = "FAPs"
$ws.Range("E3").Value = 1
$ws.Range("F3").Value = 0.3333333333333333
$ws.Range("G3").Value = 0.08805033333333334
$ws.Range("H3").Value = 0.264151
$ws.Range("I3").Value = 0.005589762818257384
$ws.Range("J3").Value = 0.005589762818257385
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 4.140873
$ws.Range("N3").Value = 12.422619
$ws.Range("O3").Value = 0.8044215857867821
$ws.Range("P3").Value = 0.8044215857867821
$ws.Range("Q3").Value = 0.364605247941
$ws.Range("R3").Value = 3.281447231469
$ws.Range("S3").Value = 0.004496525870434598
$ws.Range("T3").Value = 0.004496525870434598

# Row 4: ECs -> MuSCs
$ws.Range("A4").Value = "ECs"
$ws.Range("B4").Value = "Fgf7"
$ws.Range("C4").Value = "Fgfr2"
$ws.Range("D4").Value = "MuSCs"
$ws.Range("E4").Value = 1
$ws.Range("F4").Value = 0.3333333333333333
$ws.Range("G4").Value = 0.08805033333333334
$ws.Range("H4").Value = 0.264151
$ws.Range("I4").Value = 0.005589762818257384
$ws.Range("J4").Value = 0.005589762818257385
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 0.2847646666666667
$ws.Range("N4").Value = 0.8542940000000001
$ws.Range("O4").Value = 0.05531945672713084
$ws.Range("P4").Value = 0.05531945672713083
$ws.Range("Q4").Value = 0.02507362382155556
$ws.Range("R4").Value = 0.225662614394
$ws.Range("S4").Value = 0.0003092226423395143
$ws.Range("T4").Value = 0.0003092226423395143

# Row 5: ECs -> Resolving-Mac
$ws.Range("A5").Value = "ECs"
$ws.Range("B5").Value = "Fgf7"
$ws.Range("C5").Value = "Fgfr2"
$ws.Range("D5").Value = "Resolving-Mac"
$ws.Range("E5").Value = 1
$ws.Range("F5").Value = 0.3333333333333333
$ws.Range("G5").Value = 0.08805033333333334
$ws.Range("H5").Value = 0.264151
$ws.Range("I5").Value = 0.005589762818257384
$ws.Range("J5").Value = 0.005589762818257385
$ws.Range("K5").Value = 2
$ws.Range("L5").Value = 0.6666666666666666
$ws.Range("M5").Value = 0.007528666666666667
$ws.Range("N5").Value = 0.022586
$ws.Range("O5").Value = 0.00146254714376898
$ws.Range("P5").Value = 0.00146254714376898
$ws.Range("Q5").Value = 0.0006629016095555557
$ws.Range("R5").Value = 0.005966114486000001
$ws.Range("S5").Value = 0.000008175291644188382
$ws.Range("T5").Value = 0.000008175291644188382

# Row 6: FAPs -> ECs
$ws.Range("A6").Value = "FAPs"
$ws.Range("B6").Value = "Fgf7"
$ws.Range("C6").Value = "Fgfr2"
$ws.Range("D6").Value = "ECs"
$ws.Range("E6").Value = 3
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = 14.918018
$ws.Range("H6").Value = 44.754054
$ws.Range("I6").Value = 0.9470512964761942
$ws.Range("J6").Value = 0.9470512964761943
$ws.Range("K6").Value = 3
$ws.Range("L6").Value = 1
$ws.Range("M6").Value = 0.714474
$ws.Range("N6").Value = 2.143422
$ws.Range("O6").Value = 0.138796410342318
$ws.Range("P6").Value = 0.138796410342318
$ws.Range("Q6").Value = 10.658535992532
$ws.Range("R6").Value = 95.92682393278801
$ws.Range("S6").Value = 0.1314473203609341
$ws.Range("T6").Value = 0.1314473203609342

# Row 7: FAPs -> FAPs
$ws.Range("A7").Value = "FAPs"
$ws.Range("B7").Value = "Fgf7"
$ws.Range("C7").Value = "Fgfr2"
$ws.Range("D7").Value = "FAPs"
$ws.Range("E7").Value = 3
$ws.Range("F7").Value = 1
$ws.Range("G7").Value = 14.918018
$ws.Range("H7").Value = 44.754054
$ws.Range("I7").Value = 0.9470512964761942
$ws.Range("J7").Value = 0.9470512964761943
$ws.Range("K7").Value = 3
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 4.140873
$ws.Range("N7").Value = 12.422619
$ws.Range("O7").Value = 0.8044215857867821
$ws.Range("P7").Value = 0.8044215857867821
$ws.Range("Q7").Value = 61.77361794971401
$ws.Range("R7").Value = 555.962561547426
$ws.Range("S7").Value = 0.7618285057328081
$ws.Range("T7").Value = 0.7618285057328082

# Row 8: FAPs -> MuSCs
$ws.Range("A8").Value = "FAPs"
$ws.Range("B8").Value = "Fgf7"
$ws.Range("C8").Value = "Fgfr2"
$ws.Range("D8").Value = "MuSCs"
$ws.Range("E8").Value = 3
$ws.Range("F8").Value = 1
$ws.Range("G8").Value = 14.918018
$ws.Range("H8").Value = 44.754054
$ws.Range("I8").Value = 0.9470512964761942
$ws.Range("J8").Value = 0.9470512964761943
$ws.Range("K8").Value = 3
$ws.Range("L8").Value = 1
$ws.Range("M8").Value = 0.2847646666666667
$ws.Range("N8").Value = 0.8542940000000001
$ws.Range("O8").Value = 0.05531945672713084
$ws.Range("P8").Value = 0.05531945672713083
$ws.Range("Q8").Value = 4.248124423097335
$ws.Range("R8").Value = 38.23311980787601
$ws.Range("S8").Value = 0.05239036321378799
$ws.Range("T8").Value = 0.05239036321378798

# Row 9: FAPs -> Resolving-Mac
$ws.Range("A9").Value = "FAPs"
$ws.Range("B9").Value = "Fgf7"
$ws.Range("C9").Value = "Fgfr2"
$ws.Range("D9").Value = "Resolving-Mac"
$ws.Range("E9").Value = 3
$ws.Range("F9").Value = 1
$ws.Range("G9").Value = 14.918018
$ws.Range("H9").Value = 44.754054
$ws.Range("I9").Value = 0.9470512964761942
$ws.Range("J9").Value = 0.9470512964761943
$ws.Range("K9").Value = 2
$ws.Range("L9").Value = 0.6666666666666666
$ws.Range("M9").Value = 0.007528666666666667
$ws.Range("N9").Value = 0.022586
$ws.Range("O9").Value = 0.00146254714376898
$ws.Range("P9").Value = 0.00146254714376898
$ws.Range("Q9").Value = 0.1123127848493334
$ws.Range("R9").Value = 1.010815063644
$ws.Range("S9").Value = 0.001385107168663967
$ws.Range("T9").Value = 0.001385107168663967

# Row 10: MuSCs -> ECs
$ws.Range("A10").Value = "MuSCs"
$ws.Range("B10").Value = "Fgf7"
$ws.Range("C10").Value = "Fgfr2"
$ws.Range("D10").Value = "ECs"
$ws.Range("E10").Value = 3
$ws.Range("F10").Value = 1
$ws.Range("G10").Value = 0.7460013333333334
$ws.Range("H10").Value = 2.238004
$ws.Range("I10").Value = 0.04735894070554834
$ws.Range("J10").Value = 0.04735894070554835
$ws.Range("K10").Value = 3
$ws.Range("L10").Value = 1
$ws.Range("M10").Value = 0.714474
$ws.Range("N10").Value = 2.143422
$ws.Range("O10").Value = 0.138796410342318
$ws.Range("P10").Value = 0.138796410342318
$ws.Range("Q10").Value = 0.5329985566320001
$ws.Range("R10").Value = 4.796987009688
$ws.Range("S10").Value = 0.006573250967544796
$ws.Range("T10").Value = 0.006573250967544797

# Row 11: MuSCs -> FAPs
$ws.Range("A11").Value = "MuSCs"
$ws.Range("B11").Value = "Fgf7"
$ws.Range("C11").Value = "Fgfr2"
$ws.Range("D11").Value = "FAPs"
$ws.Range("E11").Value = 3
$ws.Range("F11").Value = 1
$ws.Range("G11").Value = 0.7460013333333334
$ws.Range("H11").Value = 2.238004
$ws.Range("I11").Value = 0.04735894070554834
$ws.Range("J11").Value = 0.04735894070554835
$ws.Range("K11").Value = 3
$ws.Range("L11").Value = 1
$ws.Range("M11").Value = 4.140873
$ws.Range("N11").Value = 12.422619
$ws.Range("O11").Value = 0.8044215857867821
$ws.Range("P11").Value = 0.8044215857867821
$ws.Range("Q11").Value = 3.089096779164
$ws.Range("R11").Value = 27.801871012476
$ws.Range("S11").Value = 0.03809655418353938
$ws.Range("T11").Value = 0.03809655418353939

# Row 12: MuSCs -> MuSCs
$ws.Range("A12").Value = "MuSCs"
$ws.Range("B12").Value = "Fgf7"
$ws.Range("C12").Value = "Fgfr2"
$ws.Range("D12").Value = "MuSCs"
$ws.Range("E12").Value = 3
$ws.Range("F12").Value = 1
$ws.Range("G12").Value = 0.7460013333333334
$ws.Range("H12").Value = 2.238004
$ws.Range("I12").Value = 0.04735894070554834
$ws.Range("J12").Value = 0.04735894070554835
$ws.Range("K12").Value = 3
$ws.Range("L12").Value = 1
$ws.Range("M12").Value = 0.2847646666666667
$ws.Range("N12").Value = 0.8542940000000001
$ws.Range("O12").Value = 0.05531945672713084
$ws.Range("P12").Value = 0.05531945672713083
$ws.Range("Q12").Value = 0.2124348210195556
$ws.Range("R12").Value = 1.911913389176
$ws.Range("S12").Value = 0.002619870871003337
$ws.Range("T12").Value = 0.002619870871003337

# Row 13: MuSCs -> Resolving-Mac
$ws.Range("A13").Value = "MuSCs"
$ws.Range("B13").Value = "Fgf7"
$ws.Range("C13").Value = "Fgfr2"
$ws.Range("D13").Value = "Resolving-Mac"
$ws.Range("E13").Value = 3
$ws.Range("F13").Value = 1
$ws.Range("G13").Value = 0.7460013333333334
$ws.Range("H13").Value = 2.238004
$ws.Range("I13").Value = 0.04735894070554834
$ws.Range("J13").Value = 0.04735894070554835
$ws.Range("K13").Value = 2
$ws.Range("L13").Value = 0.6666666666666666
$ws.Range("M13").Value = 0.007528666666666667
$ws.Range("N13").Value = 0.022586
$ws.Range("O13").Value = 0.00146254714376898
$ws.Range("P13").Value = 0.00146254714376898
$ws.Range("Q13").Value = 0.005616395371555556
$ws.Range("R13").Value = 0.05054755834400001
$ws.Range("S13").Value = 0.00006926468346082421
$ws.Range("T13").Value = 0.00006926468346082421

Write-Output "Updated Fgf7-Fgfr2: 12 rows (3 senders x 4 targets incl. Resolving-Mac)"
